$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Janice L%Atkins%NULL%1,                          Jane A H%Masoli%NULL%2,                          Jane A H%Masoli%NULL%0,                          Joao%Delgado%NULL%2,                          Joao%Delgado%NULL%0,                          Luke C%Pilling%NULL%2,                          Luke C%Pilling%NULL%0,                          Chia-Ling%Kuo%NULL%2,                          Chia-Ling%Kuo%NULL%0,                          George A%Kuchel%NULL%2,                          George A%Kuchel%NULL%0,                          David%Melzer%D.Melzer@exeter.ac.uk%1,                          Anne B%Newman%NULL%2,                          Anne B%Newman%NULL%0]'
$ws.Range('E3').Value = '[Angelo%Bianchetti%angelo.bianchetti@grupposandonato.it%1,                          R.%Rozzini%NULL%1,                          F.%Guerini%NULL%1,                          S.%Boffelli%NULL%1,                          P.%Ranieri%NULL%1,                          G.%Minelli%NULL%1,                          L.%Bianchetti%NULL%1,                          M.%Trabucchi%NULL%1]'
$ws.Range('E4').Value = '[Yan%Wan%NULL%1,                          Juan%Wu%NULL%1,                          Lihua%Ni%NULL%1,                          Qinqin%Luo%NULL%1,                          Cheng%Yuan%NULL%1,                          Fang%Fan%NULL%1,                          Hong%Liu%NULL%0,                          Changjiang%Zhang%NULL%1,                          Yuandi%Xiang%NULL%1,                          Qin%Xie%NULL%1]'
$ws.Range('E5').Value = '[Robert%De Smet%NULL%1,                          Bea%Mellaerts%NULL%1,                          Hannelore%Vandewinckele%NULL%1,                          Peter%Lybeert%NULL%1,                          Eric%Frans%NULL%1,                          Sara%Ombelet%NULL%1,                          Wim%Lemahieu%NULL%1,                          Rolf%Symons%NULL%3,                          Erwin%Ho%NULL%1,                          Johan%Frans%NULL%2,                          Annick%Smismans%NULL%1,                          Michaël R.%Laurent%NULL%1]'
$ws.Range('E6').Value = '[Marcello%Covino%marcello.covino@policlinicogemelli.it%2,                          Giuseppe%De Matteis%NULL%2,                          Giuseppe%De Matteis%NULL%0,                          Michele%Santoro%NULL%1,                          Luca%Sabia%NULL%1,                          Benedetta%Simeoni%NULL%1,                          Marcello%Candelli%NULL%1,                          Veronica%Ojetti%NULL%1,                          Francesco%Franceschi%NULL%1]'
$ws.Range('E7').Value = '[Maaike%van Gerwen%maaike.vangerwen@mountsinai.org%1,                          Mathilda%Alsen%NULL%2,                          Mathilda%Alsen%NULL%0,                          Christine%Little%NULL%1,                          Joshua%Barlow%NULL%1,                          Eric%Genden%NULL%1,                          Leonard%Naymagon%NULL%1,                          Douglas%Tremblay%NULL%1]'
$ws.Range('E8').Value = '[Jong-moon%Hwang%NULL%0,                          Ju-Hyun%Kim%NULL%1,                          Jin-Sung%Park%NULL%1,                          Min Cheol%Chang%wheel633@ynu.ac.kr%1,                          Donghwi%Park%bdome@hanmail.net%1]'
$ws.Range('E9').Value = '[Edward%Burn%NULL%1,                          Cristian%Tebé%NULL%2,                          Cristian%Tebé%NULL%0,                          Sergio%Fernandez-Bertolin%NULL%1,                          Maria%Aragon%NULL%1,                          Martina%Recalde%NULL%2,                          Martina%Recalde%NULL%0,                          Elena%Roel%NULL%1,                          Albert%Prats-Uribe%NULL%2,                          Albert%Prats-Uribe%NULL%0,                          Daniel%Prieto-Alhambra%daniel.prietoalhambra@ndorms.ox.ac.uk%2,                          Daniel%Prieto-Alhambra%daniel.prietoalhambra@ndorms.ox.ac.uk%0,                          Talita%Duarte-Salles%NULL%2,                          Talita%Duarte-Salles%NULL%0]'
$ws.Range('E10').Value = '[Paolo%Giorgi Rossi%NULL%0,                          Massimiliano%Marino%NULL%1,                          Debora%Formisano%NULL%1,                          Francesco%Venturelli%NULL%1,                          Massimo%Vicentini%NULL%2,                          Massimo%Vicentini%NULL%0,                          Roberto%Grilli%NULL%1,                          NULL%NULL%NULL%0,                          Gianluigi%Forloni%NULL%8,                          Gianluigi%Forloni%NULL%0,                          Gianluigi%Forloni%NULL%0,                          Gianluigi%Forloni%NULL%0,                          Gianluigi%Forloni%NULL%0,                          Gianluigi%Forloni%NULL%0,                          Gianluigi%Forloni%NULL%0,                          Gianluigi%Forloni%NULL%0]'
$ws.Range('E11').Value = '[Sanjana P.%Padala%NULL%1,                          Ashlyn M.%Jendro%NULL%1,                          Lillian C.%Orr%Lillian.Orr@va.gov%1]'
